$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 4 new email rows below the existing data (B3:B6)
$ws.Range("B3").Value = "User0224@yopmail.com"
$ws.Range("B4").Value = "User02880@yopmail.com"
$ws.Range("B5").Value = "User01329@yopmail.com"
$ws.Range("B6").Value = "User0721@yopmail.com"

# B2 picks up the numeric ("#,##0.00") cell style (cellXfs index 1),
# reusing the existing style instead of creating a new one.
$ws.Range("B2").NumberFormat = "#,##0.00"

# C2 collapses back onto the General style (cellXfs index 0), which is a
# duplicate of the old index-3 style that gets removed, shrinking cellXfs
# from 4 entries down to 3.
$ws.Range("C2").NumberFormat = "General"

# Leave the view scrolled/selected where the author left it when they saved.
$excel.ActiveWindow.ScrollRow = 473
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B517").Select()
